$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update cell text (row 1 headers are unchanged in content; rows 2-6 get
#    corrected/renamed labels, sources and updated release-date text per the
#    commit "Changed the data table text corrected some downloads").
# ---------------------------------------------------------------------------

$ws.Range("A2").Value = "Employment rates"
$ws.Range("B2").Value = "<a href='https://www.nomisweb.co.uk/datasets/apsnew'>Annual Population Survey</a>"
$ws.Range("C2").Value = "Jul 2021 - Jun 2022 (11/10/22)"
$ws.Range("D2").Value = "Oct 2021 - Sep 2022 (17/01/23)"

$ws.Range("A3").Value = "Employment share by occupation"
$ws.Range("B3").Value = "<a href='https://www.nomisweb.co.uk/datasets/apsnew'>Annual Population Survey</a>"
$ws.Range("C3").Value = "Jul 2021 - Jun 2022 (11/10/22)"
$ws.Range("D3").Value = "Oct 2021 - Sep 2022 (17/01/23)"

$ws.Range("A4").Value = "Online job adverts by local authority"
$ws.Range("B4").Value = "<a href='https://www.ons.gov.uk/employmentandlabourmarket/peopleinwork/employmentandemployeetypes/datasets/onlinejobadvertsbyitl1regionandlocalauthority'>ONS online job adverts</a>"
$ws.Range("C4").Value = "Jan 2022 (25/05/22)"
$ws.Range("D4").Value = "To be announced"

$ws.Range("A5").Value = "Further education and skills achievements"
$ws.Range("B5").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-tables/permalink/3960ad0f-fd8a-49bb-91d7-f3ca1181b93f'>Individualised Learner Record</a>"
$ws.Range("C5").Value = "Aug 2020 – Jul 2021 (25/11/21)"
$ws.Range("D5").Value = "Aug 2021 – Jul 2022 (Nov 22)"

$ws.Range("A6").Value = "Further education and skills achievements by sector subject area"
$ws.Range("B6").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-tables/permalink/61db0688-4ec0-4cfe-9e83-24d4ea9d078e'>Individualised Learner Record</a>"
$ws.Range("C6").Value = "Aug 2021 – Apr 2022 (provisional) (20/07/22)"
$ws.Range("D6").Value = "Aug 2021 – Jul 2022 (Nov 22)"

# ---------------------------------------------------------------------------
# 2. Apply the new small-font / left-center-aligned style to column A (rows
#    2-6) and to the two "latest" cells C5:D6 that now carry the same look.
#    (applied as two separate calls since this runtime does not support
#    multi-area Range() selectors for style assignment)
# ---------------------------------------------------------------------------

$styledRangeA = $ws.Range("A2:A6")
$styledRangeA.Font.Size = 8
$styledRangeA.Font.Color = 0
$styledRangeA.HorizontalAlignment = -4131
$styledRangeA.VerticalAlignment = -4108
$styledRangeA.ReadingOrder = 1

$styledRangeB = $ws.Range("C5:D6")
$styledRangeB.Font.Size = 8
$styledRangeB.Font.Color = 0
$styledRangeB.HorizontalAlignment = -4131
$styledRangeB.VerticalAlignment = -4108
$styledRangeB.ReadingOrder = 1

# ---------------------------------------------------------------------------
# 3. Move the active selection to B6 (matches the new saved selection state).
# ---------------------------------------------------------------------------

$ws.Range("B6").Select() | Out-Null
